# Clinical_Research_Field_List.xlsx - "Table c: extracted, cleaned, and exported into csv file"
#
# 1. On the existing "Field List" sheet:
#    - the placeholder "c" in the Table # column (C13) is marked as extracted -> "c*"
#    - the stray fill-only formatting on B16 is cleared
# 2. A new "Sheet1" worksheet is added (after "Field List") that builds up the
#    comma separated header list for the extracted/cleaned "table c" CSV export,
#    using CONCAT formulas that chain off each previous cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Field List")

# --- Field List edits -------------------------------------------------
# Table c has now been extracted -> flag it with an asterisk
$ws1.Range("C13").Value = "c*"

# Drop the leftover (fill-only, no-op) formatting that was sitting on B16
$ws1.Range("B16").ClearFormats()

# --- New "Sheet1" worksheet -------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("A1").Value = "OrgStudyId"
$ws2.Range("B1").Formula = "=A1"

$ws2.Range("A2").Value = "BriefTitle"
$ws2.Range("B2").Formula = '=CONCAT(B1,",",A2)'

$ws2.Range("A3").Value = "StartDate"
$ws2.Range("A4").Value = "CompletionDate"
$ws2.Range("A5").Value = "OverallStatus"
$ws2.Range("A6").Value = "StudyType"

# B3 carries the formula; filling it down B3:B6 creates the shared formula
# group (si="0") that auto-adjusts the relative references per row.
$ws2.Range("B3:B6").Formula = '=CONCAT(B2,",",A3)'

$ws2.Columns.Item(1).AutoFit()

# --- Selections (match where the author last clicked) -----------------
$ws2.Range("A55").Select()
$ws1.Activate()
$ws1.Range("C4").Select()
